$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.003816446910453924
$ws.Range("C2").Value = 0.0009184772660171081
$ws.Range("D2").Value = 0.001340095958043906
$ws.Range("E2").Value = [double]"1.424932756379889e-05"
$ws.Range("F2").Value = [double]"1.812596280010206e-06"
$ws.Range("G2").Value = [double]"2.297585820701277e-10"
$ws.Range("H2").Value = 0.02457530269296265
$ws.Range("I2").Value = 0.007918219843555209
$ws.Range("J2").Value = 0.01244473711807324
$ws.Range("K2").Value = 0.001120486222880649
$ws.Range("L2").Value = 0.0001035838867956089
$ws.Range("M2").Value = [double]"1.536625342293122e-08"
$ws.Range("N2").Value = 0.01025747062702243
$ws.Range("O2").Value = 0.03543484420346305
$ws.Range("P2").Value = 0.06769795001641737
$ws.Range("Q2").Value = 0.07046765132001143
$ws.Range("R2").Value = 0.006310022675344371
$ws.Range("S2").Value = [double]"5.096253608676356e-06"
$ws.Range("T2").Value = 0.02680904659869605
$ws.Range("U2").Value = 0.03759374401519867
$ws.Range("V2").Value = 0.04733750905312697
$ws.Range("W2").Value = 0.1723255406118731
$ws.Range("X2").Value = 0.02170948132728107
$ws.Range("Y2").Value = 0.0001621828073093329
$ws.Range("Z2").Value = 0.1497693723891626
$ws.Range("AA2").Value = 0.04495516917907684
$ws.Range("AB2").Value = 0.01795415700236087
$ws.Range("AC2").Value = 0.005887076843752416
$ws.Range("AD2").Value = 0.000312574583468033
$ws.Range("AE2").Value = 0.0003263981492723664
$ws.Range("AF2").Value = 0.2324312823224671
$ws.Range("B3").Value = 0.003816446910457746
$ws.Range("C3").Value = 0.0009184772660180544
$ws.Range("D3").Value = 0.001340095958044936
$ws.Range("E3").Value = [double]"1.424932756378649e-05"
$ws.Range("F3").Value = [double]"1.812596280009198e-06"
$ws.Range("G3").Value = [double]"2.297585820699941e-10"
$ws.Range("H3").Value = 0.0245753026929795
$ws.Range("I3").Value = 0.007918219843567725
$ws.Range("J3").Value = 0.01244473711808771
$ws.Range("K3").Value = 0.001120486222879716
$ws.Range("L3").Value = 0.0001035838867955509
$ws.Range("M3").Value = [double]"1.536625342292231e-08"
$ws.Range("N3").Value = 0.01025747062709203
$ws.Range("O3").Value = 0.03543484420351507
$ws.Range("P3").Value = 0.06769795001657598
$ws.Range("Q3").Value = 0.07046765131995252
$ws.Range("R3").Value = 0.006310022675339456
$ws.Range("S3").Value = [double]"5.096253608672238e-06"
$ws.Range("T3").Value = 0.02680904659844164
$ws.Range("U3").Value = 0.03759374401532316
$ws.Range("V3").Value = 0.04733750905318397
$ws.Range("W3").Value = 0.1723255406116359
$ws.Range("X3").Value = 0.02170948132725687
$ws.Range("Y3").Value = 0.0001621828073091803
$ws.Range("Z3").Value = 0.1497693723896597
$ws.Range("AA3").Value = 0.04495516917909391
$ws.Range("AB3").Value = 0.01795415700240145
$ws.Range("AC3").Value = 0.005887076843806079
$ws.Range("AD3").Value = 0.0003125745834613657
$ws.Range("AE3").Value = 0.0003263981492713747
$ws.Range("AF3").Value = 0.2324312823225808
$ws.Range("B4").Value = 0.008005961566373484
$ws.Range("C4").Value = 0.001534163580740621
$ws.Range("D4").Value = 0.002086160060275118
$ws.Range("E4").Value = [double]"1.543716967627793e-05"
$ws.Range("F4").Value = [double]"4.850602530286452e-07"
$ws.Range("G4").Value = [double]"5.865305758891648e-11"
$ws.Range("H4").Value = 0.0585410292475217
$ws.Range("I4").Value = 0.01002310542885123
$ws.Range("J4").Value = 0.01312131294096498
$ws.Range("K4").Value = 0.000488290226572279
$ws.Range("L4").Value = [double]"2.631384826434223e-05"
$ws.Range("M4").Value = [double]"3.817070681060993e-09"
$ws.Range("N4").Value = 0.199817247610114
$ws.Range("O4").Value = 0.01779258857168666
$ws.Range("P4").Value = 0.01915225772467005
$ws.Range("Q4").Value = 0.01848917457718965
$ws.Range("R4").Value = 0.001458803030729461
$ws.Range("S4").Value = [double]"1.122533030155338e-06"
$ws.Range("T4").Value = 0.06880513647487201
$ws.Range("U4").Value = 0.0005525793668040465
$ws.Range("V4").Value = 0.001300872092994791
$ws.Range("W4").Value = 0.02120831596538764
$ws.Range("X4").Value = 0.00271791910777189
$ws.Range("Y4").Value = [double]"2.468403421965914e-05"
$ws.Range("Z4").Value = 0.396355023338754
$ws.Range("AA4").Value = 0.01325960379255548
$ws.Range("AB4").Value = 0.01118016061892909
$ws.Range("AC4").Value = 0.09672739831662404
$ws.Range("AD4").Value = 0.008828186255960097
$ws.Range("AE4").Value = [double]"1.37788237315022e-05"
$ws.Range("AF4").Value = 0.02847288355727807
$ws.Range("B5").Value = [double]"1.264172661508026e-05"
$ws.Range("C5").Value = [double]"4.966349167268359e-07"
$ws.Range("D5").Value = [double]"2.888774176151704e-06"
$ws.Range("E5").Value = [double]"2.163744302740871e-05"
$ws.Range("F5").Value = [double]"1.673931365763846e-06"
$ws.Range("G5").Value = [double]"2.094331167075839e-10"
$ws.Range("H5").Value = 0.0002803292246209657
$ws.Range("I5").Value = [double]"9.573162468869423e-05"
$ws.Range("J5").Value = 0.0003220723940074238
$ws.Range("K5").Value = 0.001200293643513723
$ws.Range("L5").Value = [double]"9.429380190793072e-05"
$ws.Range("M5").Value = [double]"1.390293166185005e-08"
$ws.Range("N5").Value = 0.02606442587228198
$ws.Range("O5").Value = 0.003760404906595632
$ws.Range("P5").Value = 0.01179061100550676
$ws.Range("Q5").Value = 0.06472100759039723
$ws.Range("R5").Value = 0.005600276350312275
$ws.Range("S5").Value = [double]"4.466018189449944e-06"
$ws.Range("T5").Value = 0.06218879877540875
$ws.Range("U5").Value = 0.01492649269938281
$ws.Range("V5").Value = 0.01744377281132783
$ws.Range("W5").Value = 0.1318519803319705
$ws.Range("X5").Value = 0.01667001079874902
$ws.Range("Y5").Value = 0.000129938270854755
$ws.Range("Z5").Value = 0.3532970309932475
$ws.Range("AA5").Value = 0.04267868541630575
$ws.Range("AB5").Value = 0.02179696000275698
$ws.Range("AC5").Value = 0.04593403749381567
$ws.Range("AD5").Value = 0.001327514744275562
$ws.Range("AE5").Value = 0.0001122655115774316
$ws.Range("AF5").Value = 0.1776692430501002
$ws.Range("B6").Value = [double]"1.264172661508026e-05"
$ws.Range("C6").Value = [double]"4.966349167268413e-07"
$ws.Range("D6").Value = [double]"2.88877417615173e-06"
$ws.Range("E6").Value = [double]"2.163744302743654e-05"
$ws.Range("F6").Value = [double]"1.673931365765781e-06"
$ws.Range("G6").Value = [double]"2.094331167078235e-10"
$ws.Range("H6").Value = 0.0002803292246252394
$ws.Range("I6").Value = [double]"9.573162468869454e-05"
$ws.Range("J6").Value = 0.0003220723940051339
$ws.Range("K6").Value = 0.001200293643515381
$ws.Range("L6").Value = [double]"9.429380190805661e-05"
$ws.Range("M6").Value = [double]"1.390293166186751e-08"
$ws.Range("N6").Value = 0.02606442587236442
$ws.Range("O6").Value = 0.003760404906595632
$ws.Range("P6").Value = 0.01179061100550676
$ws.Range("Q6").Value = 0.06472100759047843
$ws.Range("R6").Value = 0.005600276350320635
$ws.Range("S6").Value = [double]"4.466018189455213e-06"
$ws.Range("T6").Value = 0.0621887987759181
$ws.Range("U6").Value = 0.01492649269938282
$ws.Range("V6").Value = 0.01744377281126039
$ws.Range("W6").Value = 0.1318519803323414
$ws.Range("X6").Value = 0.01667001079879846
$ws.Range("Y6").Value = 0.0001299382708549823
$ws.Range("Z6").Value = 0.3532970309926404
$ws.Range("AA6").Value = 0.04267868541630575
$ws.Range("AB6").Value = 0.02179696000271928
$ws.Range("AC6").Value = 0.04593403749359678
$ws.Range("AD6").Value = 0.001327514744256957
$ws.Range("AE6").Value = 0.0001122655115781079
$ws.Range("AF6").Value = 0.177669243050238
$ws.Range("B7").Value = [double]"1.26417257905111e-05"
$ws.Range("C7").Value = [double]"4.966349878496882e-07"
$ws.Range("D7").Value = [double]"2.888774375808409e-06"
$ws.Range("E7").Value = [double]"2.163744297459813e-05"
$ws.Range("F7").Value = [double]"1.673931363001298e-06"
$ws.Range("G7").Value = [double]"2.094331163662582e-10"
$ws.Range("H7").Value = 0.0002803292118162617
$ws.Range("I7").Value = [double]"9.573162777975072e-05"
$ws.Range("J7").Value = 0.0003220724004880406
$ws.Range("K7").Value = 0.001200293641049162
$ws.Range("L7").Value = [double]"9.429380174696018e-05"
$ws.Range("M7").Value = [double]"1.390293163831756e-08"
$ws.Range("N7").Value = 0.02606442544715269
$ws.Range("O7").Value = 0.003760404961165555
$ws.Range("P7").Value = 0.01179061111027229
$ws.Range("Q7").Value = 0.06472100743594703
$ws.Range("R7").Value = 0.005600276337817288
$ws.Range("S7").Value = [double]"4.466018179955741e-06"
$ws.Range("T7").Value = 0.06218879769957539
$ws.Range("U7").Value = 0.01492649282575302
$ws.Range("V7").Value = 0.01744377296496832
$ws.Range("W7").Value = 0.131851979759365
$ws.Range("X7").Value = 0.01667001073354186
$ws.Range("Y7").Value = 0.0001299382704654091
$ws.Range("Z7").Value = 0.353297033256738
$ws.Range("AA7").Value = 0.04267868547864528
$ws.Range("AB7").Value = 0.02179696011221862
$ws.Range("AC7").Value = 0.04593403789497049
$ws.Range("AD7").Value = 0.001327514772393726
$ws.Range("AE7").Value = 0.0001122655103023722
$ws.Range("AF7").Value = 0.177669243447469
$ws.Range("B8").Value = 0.008005961578339806
$ws.Range("C8").Value = 0.001534163583081304
$ws.Range("D8").Value = 0.002086160063480103
$ws.Range("E8").Value = [double]"1.543716968973624e-05"
$ws.Range("F8").Value = [double]"4.850602532738004e-07"
$ws.Range("G8").Value = [double]"5.865305761769951e-11"
$ws.Range("H8").Value = 0.05854102934994441
$ws.Range("I8").Value = 0.01002310544747319
$ws.Range("J8").Value = 0.01312131296585041
$ws.Range("K8").Value = 0.0004882902269500064
$ws.Range("L8").Value = [double]"2.63138482791274e-05"
$ws.Range("M8").Value = [double]"3.817070683163593e-09"
$ws.Range("N8").Value = 0.1998172481836647
$ws.Range("O8").Value = 0.01779258863927367
$ws.Range("P8").Value = 0.01915225780437108
$ws.Range("Q8").Value = 0.0184891745988536
$ws.Range("R8").Value = 0.001458803032287394
$ws.Range("S8").Value = [double]"1.1225330312564e-06"
$ws.Range("T8").Value = 0.06880513699596304
$ws.Range("U8").Value = 0.0005525793806084941
$ws.Range("V8").Value = 0.0013008721167946
$ws.Range("W8").Value = 0.02120831604655484
$ws.Range("X8").Value = 0.00271791911671068
$ws.Range("Y8").Value = [double]"2.468403427032042e-05"
$ws.Range("Z8").Value = 0.3963550220369417
$ws.Range("AA8").Value = 0.01325960376341723
$ws.Range("AB8").Value = 0.01118016057364305
$ws.Range("AC8").Value = 0.096727398016972
$ws.Range("AD8").Value = 0.008828186221031516
$ws.Range("AE8").Value = [double]"1.377882354586539e-05"
$ws.Range("AF8").Value = 0.0284728834405675
